$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.288864850997925
$ws.Range("B1").Value = 4.826637268066406
$ws.Range("C1").Value = 0.2777638137340546
$ws.Range("D1").Value = 0.1666877269744873
$ws.Range("E1").Value = 0.144372284412384
